# Apply weekly reshuffle of Fruta/Hortaliza data rows 2-11
# Columns updated per row: D (Fecha), J (Volumen), K (Precio minimo),
# L (Precio maximo), M (Precio promedio ponderado), P (Precio $/Kg)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2  = @{ D = 44425; J = 30; K = 13000; L = 13000; M = 13000; P = 1300 }
    3  = @{ D = 44526; J = 25; K = 9000;  L = 9000;  M = 9000;  P = 900  }
    4  = @{ D = 44523; J = 30; K = 9000;  L = 9000;  M = 9000;  P = 900  }
    5  = @{ D = 44369; J = 25; K = 8000;  L = 8000;  M = 8000;  P = 800  }
    6  = @{ D = 44348; J = 20; K = 10000; L = 10000; M = 10000; P = 1000 }
    7  = @{ D = 44530; J = 30; K = 10000; L = 10000; M = 10000; P = 1000 }
    8  = @{ D = 44463; J = 25; K = 12000; L = 12000; M = 12000; P = 1200 }
    9  = @{ D = 44473; J = 25; K = 11000; L = 11000; M = 11000; P = 1100 }
    10 = @{ D = 44525; J = 20; K = 9000;  L = 9000;  M = 9000;  P = 900  }
    11 = @{ D = 44469; J = 20; K = 12000; L = 12000; M = 12000; P = 1200 }
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Range("D$row").Value = $vals.D
    $ws.Range("J$row").Value = $vals.J
    $ws.Range("K$row").Value = $vals.K
    $ws.Range("L$row").Value = $vals.L
    $ws.Range("M$row").Value = $vals.M
    $ws.Range("P$row").Value = $vals.P
}
